# Add a new worksheet "L6" at the end of the workbook, containing
# "last six" form / goals summaries for each team.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "L6"

# Header row
$ws.Range("B1").Value = "Form"
$ws.Range("C1").Value = "Goals scored"
$ws.Range("D1").Value = "Goals conceded"
$ws.Range("E1").Value = "Total Goals"

$teams = @(
    "Airdrie Utd",
    "Clyde",
    "Cove Rangers",
    "Dumbarton",
    "East Fife",
    "Falkirk",
    "Forfar",
    "Montrose",
    "Partick",
    "Peterhead"
)

$form = @(
    "D W W L W W",
    "L L W L L W",
    "W L D D W L",
    "L W L L W L",
    "L W W L W L",
    "W L L D L L",
    "L L L L W W",
    "D L W L L W",
    "W W W W D W",
    "L L W W W L"
)

$goalsScored = @(
    "2 3 2 0 2 2",
    "0 0 3 0 1 2",
    "2 0 2 2 2 0",
    "0 1 0 1 1 0",
    "1 2 2 0 2 1",
    "2 0 0 2 0 1",
    "1 1 1 0 2 2",
    "2 0 1 1 0 2",
    "2 5 2 1 2 5",
    "0 1 1 3 3 1"
)

$goalsConceded = @(
    "2 1 0 1 1 0",
    "2 1 1 3 2 0",
    "0 1 2 2 0 2",
    "1 0 2 2 0 2",
    "2 1 1 2 1 3",
    "1 2 1 2 5 2",
    "2 3 3 1 1 1",
    "2 5 0 2 2 1",
    "0 0 0 0 2 0",
    "3 2 0 0 1 2"
)

$totalGoals = @(
    "4 4 2 1 3 2",
    "2 1 4 3 3 2",
    "2 1 4 4 2 2",
    "1 1 2 3 1 2",
    "3 3 3 2 3 4",
    "3 2 1 4 5 3",
    "3 4 4 1 3 3",
    "4 5 1 3 2 3",
    "2 5 2 1 4 5",
    "3 3 1 3 4 3"
)

# The shared-string table is built in the order cells are written, so we
# populate column by column (all of "Form" first, then "Goals scored",
# then "Goals conceded", then "Total Goals") to match the original
# authoring order, and fill column A (plain row numbers) last.
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = ($teams[$i] + "," + $form[$i])
}
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = ($teams[$i] + "," + $goalsScored[$i])
}
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = ($teams[$i] + "," + $goalsConceded[$i])
}
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = ($teams[$i] + "," + $totalGoals[$i])
}

for ($i = 0; $i -lt $teams.Length; $i++) {
    $r = $i + 2

    # Column A holds the row number (1-10), stored as text.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = [string]($i + 1)
    $cellA.ClearFormats()
}
